# Insert a new weekly price observation row into the "Poroto granado" dataset.
# This shifts the existing row 283 (and everything below it) down by one row,
# and fills the newly-opened row 283 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 283 - pushes old rows 283..344 to 284..345.
$ws.Rows.Item(283).Insert()

# Populate the new row 283 with the new observation.
$ws.Cells.Item(283, 1).Value = 9
$ws.Cells.Item(283, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(283, 3).Value = "Metropolitana"
$ws.Cells.Item(283, 4).Value = 44951
$ws.Cells.Item(283, 5).Value = 13
$ws.Cells.Item(283, 6).Value = 100112030
$ws.Cells.Item(283, 7).Value = "Poroto granado"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 70
$ws.Cells.Item(283, 11).Value = 39000
$ws.Cells.Item(283, 12).Value = 41000
$ws.Cells.Item(283, 13).Value = 40000
$ws.Cells.Item(283, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 1600
$ws.Cells.Item(283, 17).Value = 25
$ws.Cells.Item(283, 18).Value = "Hortaliza"
